$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.224.84"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.860.70"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'236.04"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4716"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("D8").Value = "'0.2898"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "'0.06569"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "'21.74"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'0.07939"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'97.62"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "1.864.26"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'5.114"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "'0.6794"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "'267.44"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "30.212.44"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "'13.67"
$ws.Range("E18").Value = "  +8.27%  "
$ws.Range("D19").Value = "'0.000007632"
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "2.106.10"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'5.250"
$ws.Range("E23").Value = "  -4.71%  "
$ws.Range("D24").Value = "'6.166"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'167.40"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'9.177"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'18.91"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "'1.946"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Value = "'0.09860"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.471"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.323"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").Value = "'4.014"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "'0.04708"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").Value = "'1.127"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").Value = "'0.7005"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "'2.710"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").Value = "'0.01873"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("D39").Value = "'2.616"
$ws.Range("E39").Value = "  +3.29%  "
$ws.Range("D40").Value = "'6.328"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").Value = "'73.64"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "'1.936"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "'0.8416"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'0.4154"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'103.35"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "'7.133"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "'940.31"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'9.152"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'34.06"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "'0.05660"
